$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the inlineStr text content of the K-column cells that now sit empty
$ws.Range("K3").ClearContents()
$ws.Range("K4").ClearContents()
$ws.Range("K5").ClearContents()
$ws.Range("K6").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("K10").ClearContents()
$ws.Range("K11").ClearContents()
$ws.Range("K12").ClearContents()
$ws.Range("K13").ClearContents()
$ws.Range("K14").ClearContents()

# Update D12 text to prepend the WU rate before the existing Diğer rate
$ws.Range("D12").Value = "WU: 0,75 USD–12 USD; Diğer: 700 TL–4.000 TL"
